$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and OKB/Filecoin row swap).
# Price/volume cells are plain text in the source sheet (e.g. thousand-dot
# separated prices, padded percent strings), so numeric-looking values are
# forced to text (temporarily via NumberFormat "@") and the cell style is
# reset to Normal afterwards so no stray formatting is introduced.

# Row 2
$ws.Range("D2").Value = '79.746.86'
$ws.Range("E2").Value = '  +5.01%  '

# Row 3
$ws.Range("D3").Value = '3.211.50'
$ws.Range("E3").Value = '  +6.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.90'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +7.14%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '639.43'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.87%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.265'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +29.77%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.606'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +10.72%  '

# Row 10
$ws.Range("D10").Value = '3.210.30'
$ws.Range("E10").Value = '  +6.70%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +41.52%  '

# Row 12
$ws.Range("E12").Value = '  +41.47%  '

# Row 13
$ws.Range("E13").Value = '  +3.67%  '

# Row 14
$ws.Range("E14").Value = '  +4.73%  '

# Row 15
$ws.Range("D15").Value = '3.800.06'
$ws.Range("E15").Value = '  +6.70%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '32.80'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +14.00%  '

# Row 17
$ws.Range("D17").Value = '79.504.38'
$ws.Range("E17").Value = '  +4.81%  '

# Row 18
$ws.Range("D18").Value = '3.205.69'
$ws.Range("E18").Value = '  +6.26%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '14.69'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +9.61%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '9.41'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.86%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +28.65%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '447.44'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +18.49%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +22.52%  '

# Row 24
$ws.Range("E24").Value = '  +13.20%  '

# Row 25
$ws.Range("D25").Value = '3.366.32'
$ws.Range("E25").Value = '  +6.39%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '77.88'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +8.15%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.96'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +12.64%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '

# Row 29
$ws.Range("E29").Value = '  +18.44%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '9.31'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +13.16%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.06%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '566.13'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +15.03%  '

# Row 33
$ws.Range("E33").Value = '  +10.71%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.157'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +30.18%  '

# Row 35
$ws.Range("E35").Value = '  +7.35%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '23.46'
$c.Style = "Normal"

# Row 37
$ws.Range("E37").Value = '  +20.43%  '

# Row 38
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("E39").Value = '  +10.46%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '163.70'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.12%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.82'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +14.14%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '20.31'
$c.Style = "Normal"

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '192.33'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.22%  '

# Row 44
$ws.Range("E44").Value = '  +0.03%  '

# Row 45
$ws.Range("E45").Value = '  +12.56%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +13.91%  '

# Row 47
$ws.Range("E47").Value = '  +4.62%  '

# Row 48
$ws.Range("E48").Value = '  +8.87%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '4.35'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +13.13%  '

# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '43.21'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.04%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '26.03'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +17.72%  '
